$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.940.07"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "2.167.79"
$ws.Range("E3").Value = "  -3.18%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "66.39"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.91%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.561"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0922"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "35.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -16.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("D15").Value = "2.491.64"
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.848"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.28%  "
$ws.Range("D18").Value = "2.180.86"
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("D19").Value = "40.849.92"
$ws.Range("E19").Value = "  -2.26%  "
$ws.Range("D20").Value = "0.0₃0937"
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("E24").Value = "  -10.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.29%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("E28").Value = "  -4.49%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.63%  "
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0733"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("E35").Value = "  -3.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0298"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.16%  "
$ws.Range("E40").Value = "  -6.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -14.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.190"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -11.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.62%  "
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0984"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.90%  "
